$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The old "总计" sheet (3rd tab) becomes the new "2022-Q1" data sheet.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"
$q1.Cells.ClearContents()

# Extend the bold/bordered header style (already present on B1:D1) across
# the new columns E1:H1 by copying the format from the existing D1 cell.
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data cells B,D,E,F,G hold text, even the numeric-looking ones (fund
# code / size / position values are stored as text in the source data).
# Force that interpretation with a leading apostrophe so Excel keeps them
# as text instead of silently parsing them into numbers. The fund-name
# column (C) is never numeric-looking so it needs no special handling.
$q1.Range("B2").Value = "'008513"
$q1.Range("C2").Value = "南方宝丰混合A"
$q1.Range("D2").Value = "'51.62"
$q1.Range("E2").Value = "'21.13"
$q1.Range("F2").Value = "'0.48"
$q1.Range("G2").Value = "'0.2478"
$q1.Range("A2").Value = 0
$q1.Range("H2").Value = 5

$q1.Range("B3").Value = "'008514"
$q1.Range("C3").Value = "南方宝丰混合C"
$q1.Range("D3").Value = "'4.72"
$q1.Range("E3").Value = "'21.13"
$q1.Range("F3").Value = "'0.48"
$q1.Range("G3").Value = "'0.0227"
$q1.Range("A3").Value = 1
$q1.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 2. Append a brand-new "总计" sheet at the end with the refreshed totals.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$total.Name = "总计"

# Reuse the header / index style from the 2022-Q1 sheet so the new sheet
# gets the same bold, centered, bordered look instead of a brand new style.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.27

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 5
$total.Range("D3").Value = 0.36

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.25

# Restore the original active tab (first sheet) so overall workbook view
# state stays the same as before the edit.
$wb.Worksheets.Item(1).Activate()
